$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New requirement: add another user record in row 4.
$ws.Range("A4").Value = "Ll"
$ws.Range("B4").Value = "Андрей"
$ws.Range("C4").Value = "Фокин"

# Password column stores text like "qwerty"/"123" (see D2/D3), so force
# text formatting before assigning a numeric-looking value - otherwise it
# would be auto-converted into a real number. Reset the cell style back
# to Normal afterwards so the cell keeps the same (default) formatting as
# the rest of the table, matching its siblings.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "123"
$ws.Range("D4").Style = "Normal"
